$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.691.63"
$ws.Range("E2").Value = '  -1.41%  '

$ws.Range("D3").Value = "'2.098.94"
$ws.Range("E3").Value = '  -0.38%  '

$ws.Range("E4").Value = '  +0.58%  '

$ws.Range("D5").Value = "'343.48"
$ws.Range("E5").Value = '  -1.99%  '

$ws.Range("E6").Value = '  +0.54%  '

$ws.Range("D7").Value = "'0.5157"
$ws.Range("E7").Value = '  -0.24%  '

$ws.Range("D8").Value = "'0.4384"
$ws.Range("E8").Value = '  -2.41%  '

$ws.Range("D9").Value = "'53.49"
$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("D10").Value = "'0.09185"
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("D11").Value = "'1.168"
$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").Value = "'24.56"
$ws.Range("E12").Value = '  -4.76%  '

$ws.Range("D13").Value = "'2.092.65"
$ws.Range("E13").Value = '  -0.55%  '

$ws.Range("D14").Value = "'6.760"
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("D15").Value = "'8.169"
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("D16").Value = "'102.08"
$ws.Range("E16").Value = '  +2.59%  '

$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("E18").Value = '  +0.59%  '

$ws.Range("D19").Value = "'21.03"
$ws.Range("E19").Value = '  -0.18%  '

$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").Value = "'6.202"
$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("D23").Value = "'29.754.40"
$ws.Range("E23").Value = '  -1.48%  '

$ws.Range("D24").Value = "'12.50"
$ws.Range("E24").Value = '  -3.21%  '

$ws.Range("D25").Value = "'2.303"
$ws.Range("E25").Value = '  -2.17%  '

$ws.Range("D26").Value = "'2.344.05"
$ws.Range("E26").Value = '  -0.49%  '

$ws.Range("E27").Value = '  -0.76%  '

$ws.Range("D28").Value = "'162.07"
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("D29").Value = "'2.501"
$ws.Range("E29").Value = '  -2.33%  '

$ws.Range("D30").Value = "'133.53"
$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("D31").Value = "'1.128"
$ws.Range("E31").Value = '  -5.11%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = "'0.1051"
$ws.Range("E32").Value = '  -1.67%  '

$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = "'1.665"
$ws.Range("E33").Value = '  +0.68%  '

$ws.Range("D34").Value = "'6.192"
$ws.Range("E34").Value = '  -1.46%  '

$ws.Range("D35").Value = "'3.946"
$ws.Range("E35").Value = '  -0.42%  '

$ws.Range("D36").Value = "'6.301"
$ws.Range("E36").Value = '  +6.02%  '

$ws.Range("D37").Value = "'10.40"
$ws.Range("E37").Value = '  +1.77%  '

$ws.Range("D38").Value = "'0.02573"
$ws.Range("E38").Value = '  -0.87%  '

$ws.Range("D39").Value = "'0.06686"
$ws.Range("E39").Value = '  -2.49%  '

$ws.Range("D40").Value = "'0.7018"
$ws.Range("E40").Value = '  +2.52%  '

$ws.Range("D41").Value = "'12.46"
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").Value = "'1.330"
$ws.Range("E42").Value = '  +5.82%  '

$ws.Range("D43").Value = "'0.2222"

$ws.Range("E44").Value = '  +6.07%  '

$ws.Range("D45").Value = "'14.31"
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").Value = "'2.312"
$ws.Range("E46").Value = '  +0.66%  '

$ws.Range("D47").Value = "'3.615"
$ws.Range("E47").Value = '  -1.62%  '

$ws.Range("D48").Value = "'0.00000000356"
$ws.Range("E48").Value = '  -1.63%  '

$ws.Range("E49").Value = '  -0.72%  '

$ws.Range("D50").Value = "'1.202"
$ws.Range("E50").Value = '  +3.17%  '

$ws.Range("D51").Value = "'81.36"
$ws.Range("E51").Value = '  -3.40%  '
